# Update the "cryptos" price/volume listing with refreshed market data.
# Column D (Price) values are forced to Text format ("@") before being
# assigned so that values like "6.40" or "57.069.15" are preserved exactly
# (trailing zeros, thousands-separator-looking dots) instead of being
# auto-converted to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.069.15"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.403.94"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.59"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +17.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.420.64"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.40"
$ws.Range("E10").Value = "  +11.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0996"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.827.53"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.015.35"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.69"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.424.09"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.30"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.26"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.405"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.523.78"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0780"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.97"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.54"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.71"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.841"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.15"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +8.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.595"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "268.30"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0530"
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.20"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0228"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.55"
$ws.Range("E49").Value = "  -6.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.871.70"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.39"
$ws.Range("E51").Value = "  -1.91%  "
